$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 1930
$ws.Range("F6").Value = 1293
$ws.Range("F7").Value = 1293
$ws.Range("F8").Value = 3
$ws.Range("F9").Value = 1602
$ws.Range("F13").Value = 1687
$ws.Range("F15").Value = 1847
$ws.Range("F18").Value = 45
$ws.Range("F19").Value = 507
$ws.Range("F20").Value = 1593
$ws.Range("F23").Value = 10
$ws.Range("F24").Value = 1094
$ws.Range("F25").Value = 2374
$ws.Range("G25").Value = 9.9
$ws.Range("F26").Value = 426
$ws.Range("F28").Value = 1014
$ws.Range("F29").Value = 4526
$ws.Range("F31").Value = 33
$ws.Range("F35").Value = 1240
$ws.Range("F36").Value = 985

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F20").Value = 183
$ws.Range("F21").Value = 10
$ws.Range("F22").Value = 10
$ws.Range("F34").Value = 469

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 2564
$ws.Range("F8").Value = 407
$ws.Range("F9").Value = 3085
$ws.Range("F10").Value = 601
$ws.Range("F11").Value = 880
$ws.Range("F12").Value = 313
$ws.Range("F13").Value = 26
$ws.Range("F14").Value = 35
$ws.Range("F15").Value = 7
$ws.Range("F16").Value = 303

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F8").Value = 3085
$ws.Range("F9").Value = 601
$ws.Range("F10").Value = 880
$ws.Range("F11").Value = 313
$ws.Range("F12").Value = 1930
$ws.Range("F14").Value = 26
$ws.Range("F15").Value = 35
$ws.Range("F16").Value = 1293
$ws.Range("F21").Value = 1687
$ws.Range("F23").Value = 1847
$ws.Range("F25").Value = 45
$ws.Range("F26").Value = 507
$ws.Range("F28").Value = 1594
$ws.Range("F30").Value = 183
$ws.Range("F31").Value = 10
$ws.Range("F34").Value = 1094
$ws.Range("F36").Value = 2374
$ws.Range("G36").Value = 9.9
$ws.Range("F37").Value = 426
$ws.Range("F39").Value = 303
$ws.Range("F41").Value = 4526
$ws.Range("F42").Value = 33
$ws.Range("F51").Value = 1240
$ws.Range("F52").Value = 985
